$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "General Description"
# ---------------------------------------------------------------------------
$ws = $wb.Sheets("General Description")
$ws.Range("B2").Value = "A test API"
$ws.Range("B7").Value = "servers description"
$ws.Range("D7").Value = ""
$ws.Range("B8").Value = "servers description"
$ws.Range("D8").Value = ""

# ---------------------------------------------------------------------------
# Sheet "Paths"
# ---------------------------------------------------------------------------
$ws = $wb.Sheets("Paths")
$ws.Range("B3").Value = "/test"
$ws.Range("C3").Value = "post"
$ws.Range("D3").Value = "Test Op"
$ws.Range("E3").Value = "Test"
$ws.Range("F3:J3").ClearContents()

# ---------------------------------------------------------------------------
# Sheet "Tags"
# ---------------------------------------------------------------------------
$ws = $wb.Sheets("Tags")
$ws.Range("A2:B2").ClearContents()

# ---------------------------------------------------------------------------
# Sheet "Parameters"
# ---------------------------------------------------------------------------
$ws = $wb.Sheets("Parameters")
$ws.Range("A2:N2").ClearContents()

# ---------------------------------------------------------------------------
# Sheet "Headers"
# ---------------------------------------------------------------------------
$ws = $wb.Sheets("Headers")
$ws.Range("A2:M2").ClearContents()

# ---------------------------------------------------------------------------
# Sheet "Schemas"
# ---------------------------------------------------------------------------
$ws = $wb.Sheets("Schemas")

# Row 3: searchCriteria/object -> id/schema
$ws.Range("A3").Value = "id"
$ws.Range("C3").Value = "ID"
$ws.Range("D3").Value = "schema"
$ws.Range("F3").Value = "TestId"

# Row 4: dateFrom -> testOperationResponse/object; clear B4,C4,F4,H4
$ws.Range("A4").Value = "testOperationResponse"
$ws.Range("B4").Value = ""
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = "object"
$ws.Range("F4").Value = ""
$ws.Range("H4").Value = ""

# Row 5: dateTo -> TestId/integer; clear B5,H5,F5; set G5,I5,N5
$ws.Range("A5").Value = "TestId"
$ws.Range("B5").Value = ""
$ws.Range("C5").Value = "ID"
$ws.Range("D5").Value = "integer"
$ws.Range("F5").Value = ""
$ws.Range("G5").Value = "int64"
$ws.Range("H5").Value = ""
$ws.Range("I5").NumberFormat = "@"
$ws.Range("I5").Value = "1"
$ws.Range("N5").NumberFormat = "@"
$ws.Range("N5").Value = "100"

# Rows 6-20: clear entirely (old testName.../TestName.../object rows removed)
$ws.Range("A6:N20").ClearContents()

# ---------------------------------------------------------------------------
# Sheet "Responses"
# ---------------------------------------------------------------------------
$ws = $wb.Sheets("Responses")
$ws.Range("A2:O2").ClearContents()
